# Update "想去人数" (F column) and one "最低票价" (G column) value for the
# sheets that hold the comic-convention data: "展览" and "全部类型".
# Both sheets mirror the same dataset, so the same edits are applied to each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new F-column ("想去人数") value
$fUpdates = @{
    3  = 504
    5  = 8482
    6  = 334
    7  = 1495
    8  = 173
    10 = 226
    11 = 241
    12 = 380
    13 = 238
    18 = 450
    19 = 1215
    20 = 155
    22 = 130
    23 = 84
    24 = 118
    25 = 64
    26 = 102
    27 = 97
}

# Map of row number -> new G-column ("最低票价") value
$gUpdates = @{
    22 = 55
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Range("F$row").Value = $fUpdates[$row]
    }

    foreach ($row in $gUpdates.Keys) {
        $ws.Range("G$row").Value = $gUpdates[$row]
    }
}
